$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New diary entry row 34
$ws.Range("A34").Value = "23 marras"

$ws.Range("B34").Value = 0.77083333333333337
$ws.Range("B34").NumberFormat = "h:mm"

$ws.Range("C34").Value = "Erottavan hypertason teoreema, kahden monikulmion leikkaustarkastelu,"
$ws.Range("C34").WrapText = $true

$ws.Rows.Item(34).RowHeight = 43.5

# Update view to match: scrolled down, selection on C34
$ws.Application.ActiveWindow.ScrollRow = 32
$ws.Range("C34").Select()
